$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data rows (2001年-2011年, rows 2-10), keep only the header row
$ws.Rows("2:10").Delete()

# Row 2: 2010年 data (previously on row 9)
$ws.Range("A2").Value = "2010年"
$ws.Range("D2").Value = 5915
$ws.Range("G2").Value = 115
$ws.Range("I2").Value = 5800

# Row 3: 2011年 data (previously on row 10)
$ws.Range("A3").Value = "2011年"
$ws.Range("D3").Value = 230
$ws.Range("G3").Value = 230

# Re-apply the "year label" formatting (bold, thin box border, centered/top)
# that column-A cells carry elsewhere on the sheet, by copying the format
# from an already-styled header cell.
$ws.Range("B1").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false
